$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current layout (before edit):
#   B1 = "meanrank" (header)
#   A2 = "16_train (GNN-MT) val delta-auprc"  B2 = 2.228205128205128
#   A3 = "16_train (RF) val delta-auprc"      B3 = 1.946153846153846
#   A4 = "16_train (PN) val delta-auprc"      B4 = 1.825641025641026
#
# Target layout (after edit): two new methods (GNN-MT-O, PN-O) are added
# to the ranking, and every meanrank value is recomputed:
#   B1 = "meanrank" (header, unchanged)
#   A2 = "16_train (GNN-MT-O) val delta-auprc"  B2 = 3.658974358974359
#   A3 = "16_train (GNN-MT) val delta-auprc"    B3 = 3.230769230769231
#   A4 = "16_train (RF) val delta-auprc"        B4 = 2.894871794871795
#   A5 = "16_train (PN) val delta-auprc"        B5 = 2.653846153846154
#   A6 = "16_train (PN-O) val delta-auprc"      B6 = 2.561538461538461

# Extend the labelled/bordered row style (currently on A2:A4) down onto the
# two newly-needed rows (A5:A6) by copying the formatting of an existing
# styled cell, so the new cells reuse the same cell style (s="1") instead
# of Excel fabricating a brand new style entry.
$ws.Range("A4").Copy()
$ws.Range("A5:A6").PasteSpecial(-4122)

# Write the new labels (row 2 is now the new GNN-MT-O entry; the rows that
# used to be 2/3/4 shift down to 3/4/5; row 6 is the new PN-O entry).
$ws.Range("A2").Value = "16_train (GNN-MT-O) val delta-auprc"
$ws.Range("A3").Value = "16_train (GNN-MT) val delta-auprc"
$ws.Range("A4").Value = "16_train (RF) val delta-auprc"
$ws.Range("A5").Value = "16_train (PN) val delta-auprc"
$ws.Range("A6").Value = "16_train (PN-O) val delta-auprc"

# Write the recomputed meanrank values.
$ws.Range("B2").Value = 3.658974358974359
$ws.Range("B3").Value = 3.230769230769231
$ws.Range("B4").Value = 2.894871794871795
$ws.Range("B5").Value = 2.653846153846154
$ws.Range("B6").Value = 2.561538461538461
